$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old column C data (superseded by column B in the new layout)
$ws.Range("C3").ClearContents()

# Insert new blank rows at target positions (descending order to keep indices stable)
$ws.Rows(29).Insert()
$ws.Rows(27).Insert()
$ws.Rows(25).Insert()
$ws.Rows(23).Insert()
$ws.Rows(21).Insert()
$ws.Rows(19).Insert()
$ws.Rows(17).Insert()
$ws.Rows(15).Insert()
$ws.Rows(13).Insert()
$ws.Rows(12).Insert()
$ws.Rows(10).Insert()
$ws.Rows(6).Insert()

# Set cell values for A and B columns
$ws.Cells.Item(1, 1).Value = "Chapter9-1"
$ws.Cells.Item(1, 2).Value = "Chapter9-1"
$ws.Cells.Item(2, 1).Value = "障害物※破壊不能"
$ws.Cells.Item(2, 2).Value = "障害物※破壊不能"
$ws.Cells.Item(3, 1).Value = "EV002"
$ws.Cells.Item(3, 2).Value = "EV002"
$ws.Cells.Item(4, 1).Value = "クラスタ誘惑"
$ws.Cells.Item(4, 2).Value = "クラスタ誘惑"
$ws.Cells.Item(5, 1).Value = "<enemy:99><CG不透明度:1,255><CG不透明度:2,100>"
$ws.Cells.Item(5, 2).Value = "<enemy:99><CG不透明度:1,255><CG不透明度:2,100>"
$ws.Cells.Item(6, 1).Value = "ーーーーー基本変更点ーーーーー"
$ws.Cells.Item(7, 1).Value = "クラスタ"
$ws.Cells.Item(7, 2).Value = "クラスタ"
$ws.Cells.Item(8, 2).Value = $ws.Cells.Item(8, 1).Value()
$ws.Cells.Item(9, 1).Value = "自動"
$ws.Cells.Item(9, 2).Value = "自動"
$ws.Cells.Item(10, 1).Value = "ーーーーーーーーアイテム生成数ーーーーーーーー"
$ws.Cells.Item(11, 1).Value = "ERS_テンプレートランダム生成 1 1 0 3 0 0"
$ws.Cells.Item(11, 2).Value = "ERS_テンプレートランダム生成 1 1 0 3 0 0"
$ws.Cells.Item(12, 1).Value = "ーーーーーーーー魔物生成数ーーーーーーーー"
$ws.Cells.Item(13, 1).Value = "夢喰い"
$ws.Cells.Item(14, 1).Value = "ERS_テンプレートランダム生成 10 1 2 3 0 0"
$ws.Cells.Item(14, 2).Value = "ERS_テンプレートランダム生成 10 1 2 3 0 0"
$ws.Cells.Item(15, 1).Value = "大夢喰い"
$ws.Cells.Item(16, 1).Value = "ERS_テンプレートランダム生成 11 1 2 3 0 0"
$ws.Cells.Item(16, 2).Value = "ERS_テンプレートランダム生成 11 1 2 3 0 0"
$ws.Cells.Item(17, 1).Value = "スイーパー"
$ws.Cells.Item(18, 1).Value = "ERS_テンプレートランダム生成 37 1 2 3 0 0"
$ws.Cells.Item(18, 2).Value = "ERS_テンプレートランダム生成 37 1 2 3 0 0"
$ws.Cells.Item(19, 1).Value = "デスコッコ"
$ws.Cells.Item(20, 1).Value = "ERS_テンプレートランダム生成 34 1 2 3 0 0"
$ws.Cells.Item(20, 2).Value = "ERS_テンプレートランダム生成 34 1 2 3 0 0"
$ws.Cells.Item(21, 1).Value = "ニャントム"
$ws.Cells.Item(22, 1).Value = "ERS_テンプレートランダム生成 32 1 2 3 0 0"
$ws.Cells.Item(22, 2).Value = "ERS_テンプレートランダム生成 32 1 2 3 0 0"
$ws.Cells.Item(23, 1).Value = "コレクター"
$ws.Cells.Item(24, 1).Value = "ERS_テンプレートランダム生成 36 1 2 3 0 0"
$ws.Cells.Item(24, 2).Value = "ERS_テンプレートランダム生成 36 1 2 3 0 0"
$ws.Cells.Item(25, 1).Value = "リーパー"
$ws.Cells.Item(26, 1).Value = "ERS_テンプレートランダム生成 46 1 2 3 0 0"
$ws.Cells.Item(26, 2).Value = "ERS_テンプレートランダム生成 46 1 2 3 0 0"
$ws.Cells.Item(27, 1).Value = "ミャウラージ"
$ws.Cells.Item(28, 1).Value = "ERS_テンプレートランダム生成 33 1 2 3 0 0"
$ws.Cells.Item(28, 2).Value = "ERS_テンプレートランダム生成 33 1 2 3 0 0"
$ws.Cells.Item(29, 1).Value = "ドゥドゥル"
$ws.Cells.Item(30, 1).Value = "ERS_テンプレートランダム生成 35 1 2 3 0 0"
$ws.Cells.Item(30, 2).Value = "ERS_テンプレートランダム生成 35 1 2 3 0 0"
$ws.Cells.Item(31, 1).Value = "ERS_テンプレートランダム生成 47 1 2 3 0 0"
$ws.Cells.Item(31, 2).Value = "ERS_テンプレートランダム生成 47 1 2 3 0 0"
